$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.560.29'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.062.86'
$ws.Range("E3").Value = '  -0.63%  '
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.53'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.620'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.80%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.57'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.45%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.389'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.61%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0792'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.79%  '
$ws.Range("E11").Value = '  -2.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.83%  '
$ws.Range("D13").Value = '2.368.98'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.763'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.42%  '
$ws.Range("D17").Value = '2.070.04'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '37.450.92'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.15'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.87%  '
$ws.Range("D21").Value = '0.0₃0832'
$ws.Range("E21").Value = '  -1.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.13%  '
$ws.Range("E24").Value = '  +0.47%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.93%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.11'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.88%  '
$ws.Range("E28").Value = '  -4.98%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.24'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.99%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.120'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.60'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.12%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0626'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.77%  '
$ws.Range("B34").Value = 'InternetComputer(DFINITY)'
$ws.Range("C34").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.65'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.49'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.29%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.85'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.28'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0223'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.23%  '
$ws.Range("D41").Value = '1.505.05'
$ws.Range("E41").Value = '  +4.08%  '
$ws.Range("B42").Value = 'InjectiveProtocol'
$ws.Range("C42").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.26'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.37%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '98.19'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.41%  '
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.91'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0953'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.18'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.03'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.66%  '
$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.96'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.63%  '
$ws.Range("D51").Value = '2.253.83'
$ws.Range("E51").Value = '  -0.63%  '
